$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "VALUE"
$ws.Range("E2").Value = 99999
$ws.Range("E3").Value = 55555

$ws.Range("E3").Select()
